$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("groups")

# Convert D2:D6 from inline string "false" to a real boolean FALSE value
$ws.Range("D2:D6").Value = $false
